# "White Hart details added"
# Adds two new pub-run log entries to Sheet1:
#   row 58 - a lockdown "Zoom" virtual catch-up (0 miles)
#   row 59 - The White Hart, Bargates (start/end at pub)
# and re-points the totals row (which moves from row 60 to row 61
# because a fresh blank spacer row is inserted above it).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new blank row above the old totals row (old row 60) so the
# previously-blank spacer rows 58/59 are freed up to hold new data,
# while a single blank spacer row still separates the data from the
# totals line. This also pushes the old footer row (63) down to 64,
# matching the new dimension A3:P64.
$ws.Rows.Item(60).Insert()

# --- Row 59 first: "The White Hart" / "Bargates" -------------------
# (entered before row 58 so the new shared-string table ends up in
# the same order as the source workbook: Bargates, Zoom, virtual,
# Lockdown virtual catch-up, Socially distanced)
$ws.Range("A59").Value = 44062
$ws.Range("C59").Value = "Bargates"
$ws.Range("B59").Value = "The White Hart"
$ws.Range("D59").Value = "start/end at pub"
$ws.Range("E59").Value = 3.37
$ws.Range("F59").Value = 0.03636574074074074
$ws.Range("G59").Formula = "=F59/E59"
$ws.Range("H59").Value = 1
$ws.Range("I59").Value = 1
$ws.Range("N59").Value = 1

# --- Row 58: lockdown Zoom virtual catch-up -------------------------
$ws.Range("A58").Value = 44013
$ws.Range("B58").Value = "Zoom"
$ws.Range("C58").Value = "Zoom"
$ws.Range("D58").Value = "virtual"
$ws.Range("F58").Value = 0.027777777777777776
$ws.Range("H58").Value = 1
$ws.Range("I58").Value = 1
$ws.Range("J58").Value = 1
$ws.Range("N58").Value = 1
$ws.Range("O58").Value = "Lockdown virtual catch-up"
$ws.Range("P58").Value = 0

# No miles were logged for the virtual catch-up, so the pre-formatted
# (but otherwise still empty) Pace cell is cleared out entirely.
$ws.Range("G58").Clear()

# --- finish off row 59 notes / total --------------------------------
$ws.Range("O59").Value = "Socially distanced"
$ws.Range("P59").Formula = "=SUM(H59:N59)*E59"

# --- fix up the totals row (now row 61) to cover the extended data --
$ws.Range("E61").Formula = "=SUM(E5:E59)"
$ws.Range("G61").Formula = "=AVERAGE(G6:G59)"
$ws.Range("H61").Formula = "=SUM(H5:H59)"
$ws.Range("I61").Formula = "=SUM(I5:I59)"
$ws.Range("J61").Formula = "=SUM(J5:J59)"
$ws.Range("K61").Formula = "=SUM(K5:K59)"
$ws.Range("L61").Formula = "=SUM(L5:L59)"
$ws.Range("M61").Formula = "=SUM(M5:M59)"
$ws.Range("N61").Formula = "=SUM(N5:N59)"
$ws.Range("P61").Formula = "=SUM(P5:P60)"

# --- restore the view: scrolled/frozen pane + active selection ------
$ws.Range("G60").Select()
